$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension grows from 4 to 7 data rows; Excel manages dimension automatically.

# Row 2: FAPs / Slitrk1 / Ptprs / ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Slitrk1"
$ws.Range("C2").Value = "Ptprs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.008164333333333334
$ws.Range("H2").Value = 0.024493
$ws.Range("I2").Value = 0.01189809984047132
$ws.Range("J2").Value = 0.01189809984047132
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.618716333333334
$ws.Range("N2").Value = 7.856149000000001
$ws.Range("O2").Value = 0.07115908183301342
$ws.Range("P2").Value = 0.07115908183301341
$ws.Range("Q2").Value = 0.02138007305077778
$ws.Range("R2").Value = 0.192420657457
$ws.Range("S2").Value = 0.0008466578602054625
$ws.Range("T2").Value = 0.0008466578602054623

# Row 3: FAPs / Slitrk1 / Ptprs / FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Slitrk1"
$ws.Range("C3").Value = "Ptprs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.008164333333333334
$ws.Range("H3").Value = 0.024493
$ws.Range("I3").Value = 0.01189809984047132
$ws.Range("J3").Value = 0.01189809984047132
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 15.503283
$ws.Range("N3").Value = 46.509849
$ws.Range("O3").Value = 0.4212748702999519
$ws.Range("P3").Value = 0.4212748702999519
$ws.Range("Q3").Value = 0.126573970173
$ws.Range("R3").Value = 1.139165731557
$ws.Range("S3").Value = 0.005012370467110434
$ws.Range("T3").Value = 0.005012370467110432

# Row 4: FAPs / Slitrk1 / Ptprs / sCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Slitrk1"
$ws.Range("C4").Value = "Ptprs"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.008164333333333334
$ws.Range("H4").Value = 0.024493
$ws.Range("I4").Value = 0.01189809984047132
$ws.Range("J4").Value = 0.01189809984047132
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 18.67887366666666
$ws.Range("N4").Value = 56.036621
$ws.Range("O4").Value = 0.5075660478670347
$ws.Range("P4").Value = 0.5075660478670347
$ws.Range("Q4").Value = 0.1525005509058889
$ws.Range("R4").Value = 1.372504958153
$ws.Range("S4").Value = 0.006039071513155424
$ws.Range("T4").Value = 0.006039071513155424

# Row 5: sCs / Slitrk1 / Ptprs / ECs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Slitrk1"
$ws.Range("C5").Value = "Ptprs"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.6780236666666667
$ws.Range("H5").Value = 2.034071
$ws.Range("I5").Value = 0.9881019001595287
$ws.Range("J5").Value = 0.9881019001595286
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.618716333333334
$ws.Range("N5").Value = 7.856149000000001
$ws.Range("O5").Value = 0.07115908183301342
$ws.Range("P5").Value = 0.07115908183301341
$ws.Range("Q5").Value = 1.775551650286556
$ws.Range("R5").Value = 15.979964852579
$ws.Range("S5").Value = 0.07031242397280796
$ws.Range("T5").Value = 0.07031242397280794

# Row 6: sCs / Slitrk1 / Ptprs / FAPs
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Slitrk1"
$ws.Range("C6").Value = "Ptprs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.6780236666666667
$ws.Range("H6").Value = 2.034071
$ws.Range("I6").Value = 0.9881019001595287
$ws.Range("J6").Value = 0.9881019001595286
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 15.503283
$ws.Range("N6").Value = 46.509849
$ws.Range("O6").Value = 0.4212748702999519
$ws.Range("P6").Value = 0.4212748702999519
$ws.Range("Q6").Value = 10.511592785031
$ws.Range("R6").Value = 94.604335065279
$ws.Range("S6").Value = 0.4162624998328415
$ws.Range("T6").Value = 0.4162624998328414

# Row 7: sCs / Slitrk1 / Ptprs / sCs
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Slitrk1"
$ws.Range("C7").Value = "Ptprs"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.6780236666666667
$ws.Range("H7").Value = 2.034071
$ws.Range("I7").Value = 0.9881019001595287
$ws.Range("J7").Value = 0.9881019001595286
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 18.67887366666666
$ws.Range("N7").Value = 56.036621
$ws.Range("O7").Value = 0.5075660478670347
$ws.Range("P7").Value = 0.5075660478670347
$ws.Range("Q7").Value = 12.66471841267678
$ws.Range("R7").Value = 113.982465714091
$ws.Range("S7").Value = 0.5015269763538793
$ws.Range("T7").Value = 0.5015269763538792
